$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before EK ("08-dec") ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns("EK:EK").Insert()

$ws1.Range("EK1").Value = "08-dec"
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 141).Value = "-"
}

# --- Sheet "Gaz": append two new rows (171, 172) ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A171").Value = "'2025-12-06"
$ws2.Range("B171").Value = 25.905
$ws2.Range("A172").Value = "'2025-12-07"
$ws2.Range("B172").Value = 25.905

# --- Sheet "CO2": append two new rows (171, 172) ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A171").Value = "'2025-12-06"
$ws3.Range("B171").Value = 81.78
$ws3.Range("A172").Value = "'2025-12-07"
$ws3.Range("B172").Value = 81.78
